$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "1.0328 at -121.90"
$ws.Range("D3").Value = "1.0154 at 117.86"

$ws.Range("C4").Value = "1.0311 at -121.98"
$ws.Range("D4").Value = "1.0134 at 117.90"

$ws.Range("B5").Value = "1.0180 at -2.55"
$ws.Range("C5").Value = "1.0401 at -121.77"
$ws.Range("D5").Value = "1.0148 at 117.83"

$ws.Range("B6").Value = "0.9940 at -3.23"
$ws.Range("C6").Value = "1.0218 at -122.22"
$ws.Range("D6").Value = "0.9960 at 117.35"

$ws.Range("B7").Value = "0.9814 at -5.62"
$ws.Range("C7").Value = "1.0592 at -122.68"
$ws.Range("D7").Value = "0.9869 at 117.30"

$ws.Range("B8").Value = "0.9779 at -5.67"
$ws.Range("D8").Value = "0.9887 at 117.47"

$ws.Range("D9").Value = "0.9905 at 117.59"

$ws.Range("B10").Value = "0.9748 at -5.87"
$ws.Range("C10").Value = "1.0616 at -122.86"
$ws.Range("D10").Value = "0.9852 at 117.31"

$ws.Range("B11").Value = "0.9814 at -5.62"
$ws.Range("C11").Value = "1.0592 at -122.68"
$ws.Range("D11").Value = "0.9869 at 117.30"

$ws.Range("B12").Value = "0.9724 at -5.59"
